# Update TPM-derived NATMI ligand/receptor metrics (Thbs1-Itga3) for rows 2-17,
# columns G:T (Ligand avg/total expression, derived specificities, Receptor
# avg/total expression, derived specificities, and Edge expression/specificity
# weights), reflecting the refreshed TPM values referenced in the commit
# "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 0..15 below correspond to worksheet rows 2..17.
# Columns 0..13 below correspond to worksheet columns G..T.
$data = New-Object 'object[,]' 16,14

$data[0,0] = 21.18599966666667
$data[0,1] = 63.557999
$data[0,2] = 0.08765141600314529
$data[0,3] = 0.08765141600314529
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 11.25749966666667
$data[0,7] = 33.772499
$data[0,8] = 0.6929800609896341
$data[0,9] = 0.6929800609896341
$data[0,10] = 238.5013841855001
$data[0,11] = 2146.512457669501
$data[0,12] = 0.0607406836076874
$data[0,13] = 0.0607406836076874

$data[1,0] = 21.18599966666667
$data[1,1] = 63.557999
$data[1,2] = 0.08765141600314529
$data[1,3] = 0.08765141600314529
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.9898276666666668
$data[1,7] = 2.969483
$data[1,8] = 0.06093101107050686
$data[1,9] = 0.06093101107050686
$data[1,10] = 20.97048861605745
$data[1,11] = 188.734397544517
$data[1,12] = 0.005340689398833248
$data[1,13] = 0.005340689398833248

$data[2,0] = 21.18599966666667
$data[2,1] = 63.557999
$data[2,2] = 0.08765141600314529
$data[2,3] = 0.08765141600314529
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 3.821582
$data[2,7] = 11.464746
$data[2,8] = 0.2352458543950409
$data[2,9] = 0.2352458543950409
$data[2,10] = 80.96403497813934
$data[2,11] = 728.6763148032541
$data[2,12] = 0.02061963224659507
$data[2,13] = 0.02061963224659507

$data[3,0] = 21.18599966666667
$data[3,1] = 63.557999
$data[3,2] = 0.08765141600314529
$data[3,3] = 0.08765141600314529
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.1761463333333333
$data[3,7] = 0.528439
$data[3,8] = 0.01084307354481826
$data[3,9] = 0.01084307354481827
$data[3,10] = 3.731836159284555
$data[3,11] = 33.586525433561
$data[3,12] = 0.0009504107500295647
$data[3,13] = 0.000950410750029565

$data[4,0] = 121.8208923333333
$data[4,1] = 365.462677
$data[4,2] = 0.5040014103551328
$data[4,3] = 0.5040014103551328
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 11.25749966666667
$data[4,7] = 33.772499
$data[4,8] = 0.6929800609896341
$data[4,9] = 0.6929800609896341
$data[4,10] = 1371.398654835536
$data[4,11] = 12342.58789351982
$data[4,12] = 0.3492629280867615
$data[4,13] = 0.3492629280867615

$data[5,0] = 121.8208923333333
$data[5,1] = 365.462677
$data[5,2] = 0.5040014103551328
$data[5,3] = 0.5040014103551328
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.9898276666666668
$data[5,7] = 2.969483
$data[5,8] = 0.06093101107050686
$data[5,9] = 0.06093101107050686
$data[5,10] = 120.5816896095546
$data[5,11] = 1085.235206485991
$data[5,12] = 0.03070931551389967
$data[5,13] = 0.03070931551389967

$data[6,0] = 121.8208923333333
$data[6,1] = 365.462677
$data[6,2] = 0.5040014103551328
$data[6,3] = 0.5040014103551328
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 3.821582
$data[6,7] = 11.464746
$data[6,8] = 0.2352458543950409
$data[6,9] = 0.2352458543950409
$data[6,10] = 465.5485293650047
$data[6,11] = 4189.936764285042
$data[6,12] = 0.1185642423952988
$data[6,13] = 0.1185642423952988

$data[7,0] = 121.8208923333333
$data[7,1] = 365.462677
$data[7,2] = 0.5040014103551328
$data[7,3] = 0.5040014103551328
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.1761463333333333
$data[7,7] = 0.528439
$data[7,8] = 0.01084307354481826
$data[7,9] = 0.01084307354481827
$data[7,10] = 21.45830350791145
$data[7,11] = 193.124731571203
$data[7,12] = 0.005464924359172833
$data[7,13] = 0.005464924359172835

$data[8,0] = 37.20718233333333
$data[8,1] = 111.621547
$data[8,2] = 0.1539347809079331
$data[8,3] = 0.1539347809079331
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 11.25749966666667
$data[8,7] = 33.772499
$data[8,8] = 0.6929800609896341
$data[8,9] = 0.6929800609896341
$data[8,10] = 418.8598427151058
$data[8,11] = 3769.738584435952
$data[8,12] = 0.1066737338620054
$data[8,13] = 0.1066737338620055

$data[9,0] = 37.20718233333333
$data[9,1] = 111.621547
$data[9,2] = 0.1539347809079331
$data[9,3] = 0.1539347809079331
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.9898276666666668
$data[9,7] = 2.969483
$data[9,8] = 0.06093101107050686
$data[9,9] = 0.06093101107050686
$data[9,10] = 36.82869847224455
$data[9,11] = 331.458286250201
$data[9,12] = 0.00937940183963732
$data[9,13] = 0.009379401839637322

$data[10,0] = 37.20718233333333
$data[10,1] = 111.621547
$data[10,2] = 0.1539347809079331
$data[10,3] = 0.1539347809079331
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 3.821582
$data[10,7] = 11.464746
$data[10,8] = 0.2352458543950409
$data[10,9] = 0.2352458543950409
$data[10,10] = 142.1902982757846
$data[10,11] = 1279.712684482062
$data[10,12] = 0.03621251905580015
$data[10,13] = 0.03621251905580016

$data[11,0] = 37.20718233333333
$data[11,1] = 111.621547
$data[11,2] = 0.1539347809079331
$data[11,3] = 0.1539347809079331
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.1761463333333333
$data[11,7] = 0.528439
$data[11,8] = 0.01084307354481826
$data[11,9] = 0.01084307354481827
$data[11,10] = 6.553908741681443
$data[11,11] = 58.98517867513299
$data[11,12] = 0.001669126150490205
$data[11,13] = 0.001669126150490206

$data[12,0] = 61.49336899999999
$data[12,1] = 184.480107
$data[12,2] = 0.2544123927337887
$data[12,3] = 0.2544123927337887
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 11.25749966666667
$data[12,7] = 33.772499
$data[12,8] = 0.6929800609896341
$data[12,9] = 0.6929800609896341
$data[12,10] = 692.2615810197102
$data[12,11] = 6230.354229177391
$data[12,12] = 0.1763027154331797
$data[12,13] = 0.1763027154331797

$data[13,0] = 61.49336899999999
$data[13,1] = 184.480107
$data[13,2] = 0.2544123927337887
$data[13,3] = 0.2544123927337887
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 0.9898276666666668
$data[13,7] = 2.969483
$data[13,8] = 0.06093101107050686
$data[13,9] = 0.06093101107050686
$data[13,10] = 60.86783795274233
$data[13,11] = 547.810541574681
$data[13,12] = 0.01550160431813662
$data[13,13] = 0.01550160431813662

$data[14,0] = 61.49336899999999
$data[14,1] = 184.480107
$data[14,2] = 0.2544123927337887
$data[14,3] = 0.2544123927337887
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 3.821582
$data[14,7] = 11.464746
$data[14,8] = 0.2352458543950409
$data[14,9] = 0.2352458543950409
$data[14,10] = 235.001952089758
$data[14,11] = 2115.017568807822
$data[14,12] = 0.05984946069734682
$data[14,13] = 0.05984946069734683

$data[15,0] = 61.49336899999999
$data[15,1] = 184.480107
$data[15,2] = 0.2544123927337887
$data[15,3] = 0.2544123927337887
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 0.1761463333333333
$data[15,7] = 0.528439
$data[15,8] = 0.01084307354481826
$data[15,9] = 0.01084307354481827
$data[15,10] = 10.83183147366366
$data[15,11] = 97.48648326297298
$data[15,12] = 0.002758612285125658
$data[15,13] = 0.002758612285125659

$ws.Range("G2:T17").Value2 = $data
